$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "L1cam"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 19.72083766666667
$ws.Range("H2").Value = 59.162513
$ws.Range("I2").Value = 0.8016210077351786
$ws.Range("J2").Value = 0.8016210077351787
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 25.77723838635111
$ws.Range("R2").Value = 231.99514547716
$ws.Range("S2").Value = 0.0100704693624907
$ws.Range("T2").Value = 0.0100704693624907

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "L1cam"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 19.72083766666667
$ws.Range("H3").Value = 59.162513
$ws.Range("I3").Value = 0.8016210077351786
$ws.Range("J3").Value = 0.8016210077351787
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 1582.128596484527
$ws.Range("R3").Value = 14239.15736836074
$ws.Range("S3").Value = 0.6180948214706404
$ws.Range("T3").Value = 0.6180948214706405

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "L1cam"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.72083766666667
$ws.Range("H4").Value = 59.162513
$ws.Range("I4").Value = 0.8016210077351786
$ws.Range("J4").Value = 0.8016210077351787
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.541576
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 443.9921520156098
$ws.Range("R4").Value = 3995.929368140489
$ws.Range("S4").Value = 0.1734557169020475
$ws.Range("T4").Value = 0.1734557169020475

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "L1cam"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.099159
$ws.Range("H5").Value = 0.297477
$ws.Range("I5").Value = 0.004030657259573097
$ws.Range("J5").Value = 0.004030657259573097
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 0.12961138996
$ws.Range("R5").Value = 1.16650250964
$ws.Range("S5").Value = 0.00005063566205420727
$ws.Range("T5").Value = 0.00005063566205420728

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "L1cam"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.099159
$ws.Range("H6").Value = 0.297477
$ws.Range("I6").Value = 0.004030657259573097
$ws.Range("J6").Value = 0.004030657259573097
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 7.955153434682998
$ws.Range("R6").Value = 71.59638091214698
$ws.Range("S6").Value = 0.00310786313635159
$ws.Range("T6").Value = 0.00310786313635159

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "L1cam"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.099159
$ws.Range("H7").Value = 0.297477
$ws.Range("I7").Value = 0.004030657259573097
$ws.Range("J7").Value = 0.004030657259573097
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.541576
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 2.232451711528
$ws.Range("R7").Value = 20.092065403752
$ws.Range("S7").Value = 0.0008721584611672998
$ws.Range("T7").Value = 0.0008721584611673

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "L1cam"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.781202
$ws.Range("H8").Value = 14.343606
$ws.Range("I8").Value = 0.1943483350052483
$ws.Range("J8").Value = 0.1943483350052483
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.307106666666667
$ws.Range("N8").Value = 3.92132
$ws.Range("O8").Value = 0.01256263154946851
$ws.Range("P8").Value = 0.01256263154946851
$ws.Range("Q8").Value = 6.249541008880001
$ws.Range("R8").Value = 56.24586907992001
$ws.Range("S8").Value = 0.002441526524923607
$ws.Range("T8").Value = 0.002441526524923607

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "L1cam"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.781202
$ws.Range("H9").Value = 14.343606
$ws.Range("I9").Value = 0.1943483350052483
$ws.Range("J9").Value = 0.1943483350052483
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.77105616682495
$ws.Range("P9").Value = 0.77105616682495
$ws.Range("Q9").Value = 383.5778447968739
$ws.Range("R9").Value = 3452.200603171866
$ws.Range("S9").Value = 0.149853482217958
$ws.Range("T9").Value = 0.149853482217958

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "L1cam"
$ws.Range("C10").Value = "Egfr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.781202
$ws.Range("H10").Value = 14.343606
$ws.Range("I10").Value = 0.1943483350052483
$ws.Range("J10").Value = 0.1943483350052483
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.51385866666667
$ws.Range("N10").Value = 67.541576
$ws.Range("O10").Value = 0.2163812016255815
$ws.Range("P10").Value = 0.2163812016255815
$ws.Range("Q10").Value = 107.643306084784
$ws.Range("R10").Value = 968.7897547630562
$ws.Range("S10").Value = 0.04205332626236668
$ws.Range("T10").Value = 0.04205332626236669
